{"js": "// Remove the trailing clause \", until Smartcash reaches a considerable\n// market cap\" from the SmartCash-mining paragraph (keep the final period).\nconst oldFragment =\n  \", until Smartcash reaches a considerable market cap.\";\nconst newFragment = \".\";\n\nconst results = context.document.body.search(oldFragment, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newFragment, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// The \"exchanges\" bookmark around the EXCHANGES heading stays in place;\n// re-anchor it so it is preserved (by name) after the edit/save pass.\nconst bookmarkRange = context.document.getBookmarkRangeOrNullObject(\"exchanges\");\nawait context.sync();\nif (!bookmarkRange.isNullObject) {\n  context.document.deleteBookmark(\"exchanges\");\n  bookmarkRange.insertBookmark(\"exchanges\");\n  await context.sync();\n}\n", "ps1": "# Remove the trailing clause \", until Smartcash reaches a considerable\n# market cap\" from the SmartCash-mining paragraph (keep the final period).\n$d = $word.ActiveDocument\n\n$oldFragment = \", until Smartcash reaches a considerable market cap.\"\n$newFragment = \".\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n#   ReplaceWith, Replace)\n$find.Execute($oldFragment, $true, $false, $false, $false, $false, $true, 1, $false, $newFragment, 2) | Out-Null\n\n# The \"exchanges\" bookmark around the EXCHANGES heading stays in place;\n# re-anchor it so it is preserved (by name) after the edit/save pass.\nif ($d.Bookmarks.Exists(\"exchanges\")) {\n    $bookmarkRange = $d.Bookmarks.Item(\"exchanges\").Range\n    $d.Bookmarks.Item(\"exchanges\").Delete()\n    $d.Bookmarks.Add(\"exchanges\", $bookmarkRange) | Out-Null\n}\n"}
